# Update countries & provincias Spain
# - Refresh COVID-19 stats for several countries (including the two
#   country-pairs that swapped rank/row order: Costa Rica/Uzbekistan
#   and Siria/Mali).
# - Update the "last updated" footer timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($Row, $Country, $CasosTotales, $NuevosCasos, $CasosActivos, $Recuperados, $CasosCriticos, $MuertesHoy, $Muertes) {
    $ws.Cells.Item($Row, 1).Value = $Country
    $ws.Cells.Item($Row, 2).Value = $CasosTotales
    $ws.Cells.Item($Row, 3).Value = $NuevosCasos
    $ws.Cells.Item($Row, 4).Value = $CasosActivos
    $ws.Cells.Item($Row, 5).Value = $Recuperados
    $ws.Cells.Item($Row, 6).Value = $CasosCriticos
    $ws.Cells.Item($Row, 7).Value = $MuertesHoy
    $ws.Cells.Item($Row, 8).Value = $Muertes
}

# Estados Unidos (row 4)
Set-Row 4 "Estados Unidos" 6247079 31487 3478814 2579632 0 897 188633

# Sudafrica (row 9)
Set-Row 9 "Sudafrica" 628259 1218 549993 64003 0 114 14263

# Alemania (row 23)
Set-Row 23 "Alemania" 245984 1192 219900 16703 0 10 9381

# Costa Rica overtakes Uzbekistan -> rows 62/63 swap identity.
# Row 62: now Costa Rica (updated stats)
Set-Row 62 "Costa Rica" 42184 897 16270 25471 0 7 443

# Row 63: now Uzbekistan (stats previously held by row 62, unchanged)
Set-Row 63 "Uzbekistan" 42127 234 39538 2266 0 3 323

# Kenia (row 68)
Set-Row 68 "Kenia" 34315 114 20211 13527 0 0 577

# Namibia (row 103)
Set-Row 103 "Namibia" 7692 142 3379 4232 0 6 81

# Cabo Verde (row 122)
Set-Row 122 "Cabo Verde" 3970 86 3423 507 0 0 40

# Siria overtakes Mali -> rows 132/133 swap identity.
# Row 132: now Siria (updated stats)
Set-Row 132 "Siria" 2830 65 646 2068 0 4 116

# Row 133: now Mali (updated stats)
Set-Row 133 "Mali" 2777 1 2178 473 0 0 126

# Angola (row 134)
Set-Row 134 "Angola" 2729 75 1084 1536 0 1 109

# Footer timestamp update
$ws.Range("A1").Value = "Datos actualizados a 1 de Septiembre de 2020 a las 22:57"
